# Refresh the cryptocurrency price/volume table on Sheet1 with the
# latest values (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.753.80'
$ws.Range('E2').Value = '  -0.16%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.873.97'
$ws.Range('E3').Value = '  -0.64%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.010'
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '336.25'
$ws.Range('E5').Value = '  +0.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.010'
$ws.Range('E6').Value = '  +0.32%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4685'
$ws.Range('E7').Value = '  -0.56%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3941'
$ws.Range('E8').Value = '  +0.29%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.64'
$ws.Range('E9').Value = '  -4.22%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08026'
$ws.Range('E10').Value = '  -0.53%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.009'
$ws.Range('E11').Value = '  -1.68%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.98'
$ws.Range('E12').Value = '  -0.82%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.877.18'
$ws.Range('E13').Value = '  -0.32%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.015'
$ws.Range('E14').Value = '  +0.62%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.289'
$ws.Range('E15').Value = '  +2.29%  '
$ws.Range('E16').Value = '  +0.32%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '89.06'
$ws.Range('E17').Value = '  +2.00%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06754'
$ws.Range('E18').Value = '  +0.13%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.00001047'
$ws.Range('E19').Value = '  -0.23%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.33'
$ws.Range('E20').Value = '  -0.25%  '
$ws.Range('E21').Value = '  +0.35%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '27.704.49'
$ws.Range('E22').Value = '  -0.43%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.501'
$ws.Range('E23').Value = '  -0.27%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.98'
$ws.Range('E24').Value = '  -0.24%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.314'
$ws.Range('E25').Value = '  -1.13%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.098.45'
$ws.Range('E26').Value = '  -0.37%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '159.61'
$ws.Range('E27').Value = '  +0.47%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.82'
$ws.Range('E28').Value = '  -1.78%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.147'
$ws.Range('E29').Value = '  +2.00%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.466'
$ws.Range('E30').Value = '  -2.05%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '122.17'
$ws.Range('E31').Value = '  +0.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9887'
$ws.Range('E32').Value = '  +0.62%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09545'
$ws.Range('E33').Value = '  +0.57%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.642'
$ws.Range('E34').Value = '  +0.67%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.340'
$ws.Range('E35').Value = '  -0.36%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.347'
$ws.Range('E36').Value = '  -7.14%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06075'
$ws.Range('E37').Value = '  -1.41%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02242'
$ws.Range('E38').Value = '  -1.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.202'
$ws.Range('E39').Value = '  -1.41%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.329'
$ws.Range('E40').Value = '  +3.11%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.009'
$ws.Range('E41').Value = '  +0.27%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5982'
$ws.Range('E42').Value = '  -0.40%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1895'
$ws.Range('E43').Value = '  +0.12%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.47'
$ws.Range('E44').Value = '  +1.37%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5683'
$ws.Range('E45').Value = '  -0.59%  '
$ws.Range('B46').Value = 'WEMIXTOKEN'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.245'
$ws.Range('E46').Value = '  -1.23%  '
$ws.Range('E47').Value = '  -0.68%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.930'
$ws.Range('E48').Value = '  -0.79%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06775'
$ws.Range('E49').Value = '  -2.03%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '112.43'
$ws.Range('E50').Value = '  -1.54%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.028'
$ws.Range('E51').Value = '  -10.79%  '
